{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Replaces the \"history/governance\" essay content with the\n// \"government and politics\" essay content, per the commit diff:\n//  - Title, author name, and email are swapped for new values.\n//  - The body paragraph's three opening sentences are replaced with new\n//    sentences, and the remaining (now removed) sentences/line-breaks are\n//    dropped.\n//  - The \"Summary\" section text is replaced with new summary text\n//    (also removing the embedded page-break run).\n//  - A new empty paragraph is appended at the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Title\nparagraphs.items[0].insertText(\n  \"A Glimpse into the Art of Government: An Exploration of Civics and Politics\",\n  \"Replace\"\n);\n\n// 2) Author name\nparagraphs.items[1].insertText(\"Clara Bennett\", \"Replace\");\n\n// 3) Email address - rebuild it run-by-run to mirror the original\n//    run layout: \"clara\" + \".\" + \"bennett88@institute\" + \".\" + \"edu\"\nconst emailParagraph = paragraphs.items[2];\nconst emailRange = emailParagraph.getRange();\nemailRange.insertText(\"clara.bennett88@institute\", \"Replace\");\nawait context.sync();\n\nconst instituteHits = body.search(\"clara.bennett88@institute\", { matchCase: true });\ninstituteHits.load(\"items\");\nawait context.sync();\ninstituteHits.items[0].insertText(\".edu\", \"End\");\n\n// 4) Main body paragraph (index 4): collapse every run/line-break in the\n//    paragraph down to the three new sentences.\nconst newBodyText =\n  \"The exploration of government and politics provides a lens through which \" +\n  \"we can examine the interplay of power dynamics, decision-making processes, \" +\n  \"and the quest for justice. It encourages us to think critically, to \" +\n  \"challenge assumptions, and to recognize the interconnections between our \" +\n  \"actions and their broader implications. As we navigate the complexities of \" +\n  \"governance and political engagement, we gain a deeper appreciation for the \" +\n  \"rights and responsibilities that come with being a citizen, and we embrace \" +\n  \"the opportunity to contribute to a better future for ourselves and for \" +\n  \"generations to come.\";\nparagraphs.items[4].insertText(newBodyText, \"Replace\");\n\n// 5) \"Summary\" heading (index 5) is unchanged.\n\n// 6) Summary body paragraph (index 6): replace with the new summary text.\nconst newSummaryText =\n  \"In conclusion, government and politics are fundamental pillars of human \" +\n  \"society, shaping the structures, processes, and relationships that define \" +\n  \"how we live together. The study of government and politics provides a \" +\n  \"critical lens through which we can examine the intricacies of governance, \" +\n  \"the interplay of power, and the quest for a just and equitable society. It \" +\n  \"equips us with the knowledge, skills, and values necessary to navigate the \" +\n  \"political landscape, to participate effectively in the decision-making \" +\n  \"process, and to work towards a better future for all.\";\nparagraphs.items[6].insertText(newSummaryText, \"Replace\");\n\n// 7) Append a new, empty paragraph at the very end of the document body.\nbody.insertParagraph(\"\", \"End\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d (ActiveDocument) are pre-seeded by the host.\n#\n# Replaces the \"history/governance\" essay content with the\n# \"government and politics\" essay content, per the commit diff:\n#  - Title, author name, and email are swapped for new values.\n#  - The body paragraph's three opening sentences are replaced with new\n#    sentences, and the remaining (now removed) sentences/line-breaks are\n#    dropped.\n#  - The \"Summary\" section text is replaced with new summary text\n#    (also removing the embedded page-break run).\n#  - A new empty paragraph is appended at the end of the document.\n\n$d = $word.ActiveDocument\n\n# wdReplaceOne = 1 (unused here); we always pass ReplaceAll below.\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# 1) Title (paragraph 1)\n$d.Paragraphs(1).Range.Text = \"A Glimpse into the Art of Government: An Exploration of Civics and Politics\"\n\n# 2) Author name (paragraph 2)\n$d.Paragraphs(2).Range.Text = \"Clara Bennett\"\n\n# 3) Email address (paragraph 3) - rebuild run-by-run so the final text\n#    reads \"clara.bennett88@institute.edu\":\n#      \"marcuswalton56@abromail\" -> \"clara\"\n#      \".\"                        (unchanged middle run)\n#      \"net\"                     -> \"bennett88@institute\"\n#      + new run \".\"\n#      + new run \"edu\"\n$rng = $d.Paragraphs(3).Range\n$rng.Find.Execute(\"marcuswalton56@abromail\", $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, \"clara\", $wdReplaceAll) | Out-Null\n\n$rng = $d.Paragraphs(3).Range\n$rng.Find.Execute(\"net\", $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, \"bennett88@institute\", $wdReplaceAll) | Out-Null\n\n$emailEnd = $d.Paragraphs(3).Range\n$emailEnd.MoveEnd(1, -1) | Out-Null\n$emailEnd.InsertAfter(\".\")\n\n$emailEnd2 = $d.Paragraphs(3).Range\n$emailEnd2.MoveEnd(1, -1) | Out-Null\n$emailEnd2.InsertAfter(\"edu\")\n\n# 4) Main body paragraph (paragraph 5): replace the first two sentences,\n#    then delete/replace everything from \"By delving\" through the end of\n#    the paragraph (this also removes the two blank-line <w:br/> pairs)\n#    with the new third sentence.\n$rng = $d.Paragraphs(5).Range\n$rng.Find.Execute( `\n  \"History, like a murmuring river, whispers tales of triumphs and tribulations, inviting us to glean wisdom from the annals of time\", `\n  $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, `\n  \"The exploration of government and politics provides a lens through which we can examine the interplay of power dynamics, decision-making processes, and the quest for justice\", `\n  $wdReplaceAll) | Out-Null\n\n$rng = $d.Paragraphs(5).Range\n$rng.Find.Execute( `\n  \" It is a tapestry woven with threads of human experience, where patterns emerge, offering guidance for the complexities of modern governance\", `\n  $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, `\n  \" It encourages us to think critically, to challenge assumptions, and to recognize the interconnections between our actions and their broader implications\", `\n  $wdReplaceAll) | Out-Null\n\n$rng = $d.Paragraphs(5).Range\n$rng.Find.Execute( `\n  \" By delving*welfare\", `\n  $false, $false, $true, $false, $false, $true, $wdFindContinue, $false, `\n  \" As we navigate the complexities of governance and political engagement, we gain a deeper appreciation for the rights and responsibilities that come with being a citizen, and we embrace the opportunity to contribute to a better future for ourselves and for generations to come\", `\n  $wdReplaceAll) | Out-Null\n\n# 5) \"Summary\" heading (paragraph 6) is unchanged.\n\n# 6) Summary body paragraph (paragraph 7): replace each sentence with the\n#    new summary text; the third replacement also swallows the\n#    lastRenderedPageBreak run that previously split \"History serves as a\n#    constant \" / \"reminder of ... road ahead\" into two runs.\n$rng = $d.Paragraphs(7).Range\n$rng.Find.Execute( `\n  \"History, as a mirror to the present, holds a wealth of lessons for modern governance\", `\n  $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, `\n  \"In conclusion, government and politics are fundamental pillars of human society, shaping the structures, processes, and relationships that define how we live together\", `\n  $wdReplaceAll) | Out-Null\n\n$rng = $d.Paragraphs(7).Range\n$rng.Find.Execute( `\n  \" By studying the echoes of the past - the triumphs and tribulations of civilizations, the struggles for justice and equality, and the intricacies of diplomacy and statecraft - leaders can gain insights into the challenges they face and the paths they must tread\", `\n  $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, `\n  \" The study of government and politics provides a critical lens through which we can examine the intricacies of governance, the interplay of power, and the quest for a just and equitable society\", `\n  $wdReplaceAll) | Out-Null\n\n$rng = $d.Paragraphs(7).Range\n$rng.Find.Execute( `\n  \" History serves as a constant *road ahead\", `\n  $false, $false, $true, $false, $false, $true, $wdFindContinue, $false, `\n  \" It equips us with the knowledge, skills, and values necessary to navigate the political landscape, to participate effectively in the decision-making process, and to work towards a better future for all\", `\n  $wdReplaceAll) | Out-Null\n\n# 7) Append a new, empty paragraph at the very end of the document body.\n$d.Content.InsertParagraphAfter()\n"}
